$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Category" value for the row-5 sample product
$ws.Range("I5").Value = "Zarmalang"

# Simplify the dummy-product row: drop the Display Name / Product Id /
# Brief Description detail cells, keeping only Operation, Product Id and
# Category, to demonstrate "skipping details to just publish".
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("F6").ClearContents()

# Update the active selection to match where the author left off editing.
$ws.Range("F6").Select()
